# binary-classification-results.xlsx
# "Updated top 5 models results for Thunderbird Comments data"
#
# The Thunderbird sheet's F1/Accuracy/Precision/Recall columns (C:F) are
# refreshed with new scores for the top-5 models across every configuration
# block (rows 3-61, skipping the blank separator rows). Row 41 only has its
# Recall (F) value touched. The workbook is also left with the Thunderbird
# tab active/selected, matching the reviewer's final position.

$wb = $excel.ActiveWorkbook
$wsLucene = $wb.Worksheets.Item("Lucene")
$wsUbuntu = $wb.Worksheets.Item("Ubuntu")
$wsThunderbird = $wb.Worksheets.Item("Thunderbird")

function Set-ThunderbirdRow($row, $f1, $accuracy, $precision, $recall) {
    $wsThunderbird.Cells.Item($row, 3).Value = $f1
    $wsThunderbird.Cells.Item($row, 4).Value = $accuracy
    $wsThunderbird.Cells.Item($row, 5).Value = $precision
    $wsThunderbird.Cells.Item($row, 6).Value = $recall
}

# Count Vectorizer
Set-ThunderbirdRow 3 67.6 69.96 71.33 64.78
Set-ThunderbirdRow 4 69.78 66.07 61.68 81.58
Set-ThunderbirdRow 5 19.09 54.7 75.33 11.21
Set-ThunderbirdRow 6 62.35 64.84 65.17 60.4
Set-ThunderbirdRow 7 58.54 64.64 68.92 51.89

# CV + tfidf
Set-ThunderbirdRow 9 68.18 68.59 67.47 70.03
Set-ThunderbirdRow 10 68.31 66.75 64.22 74.4
Set-ThunderbirdRow 11 5.47 48.63 3.76 10
Set-ThunderbirdRow 12 61.15 62.73 62.54 60.57
Set-ThunderbirdRow 13 61.97 65.71 68.02 57.73

# CV + tfidf + ngram(1)
Set-ThunderbirdRow 15 68.18 68.59 67.47 70.03
Set-ThunderbirdRow 16 68.31 66.75 64.22 74.4
Set-ThunderbirdRow 17 5.47 48.63 3.76 10
Set-ThunderbirdRow 18 61.15 62.73 62.54 60.57
Set-ThunderbirdRow 19 61.97 65.71 68.02 57.73

# CV + tfidf + ngram(2)
Set-ThunderbirdRow 21 72.15 70.35 66.66 79.11
Set-ThunderbirdRow 22 70.19 67.52 63.94 79.19
Set-ThunderbirdRow 23 5.47 48.63 3.76 10
Set-ThunderbirdRow 24 56.47 60.85 62.14 52.47
Set-ThunderbirdRow 25 57.82 64.9 70.24 49.86

# CV + tfidf + ngram(3)
Set-ThunderbirdRow 27 73.28 70.44 65.64 83.97
Set-ThunderbirdRow 28 69.99 67.36 63.91 78.85
Set-ThunderbirdRow 29 5.47 48.63 3.76 10
Set-ThunderbirdRow 30 52.3 58.75 60.07 46.66
Set-ThunderbirdRow 31 52.65 62.7 69.97 43.07

# CV + tfidf + ngram(1) + stopwords
Set-ThunderbirdRow 33 62.12 64.22 64.36 61.38
Set-ThunderbirdRow 34 63.45 63.74 62.61 65.53
Set-ThunderbirdRow 35 5.47 48.63 3.76 10
Set-ThunderbirdRow 36 57.44 60.56 60.97 55.4
Set-ThunderbirdRow 37 53.88 61.73 66.18 46.27

# CV + tfidf + ngram(1) + stopwords + lemmatization (only Recall for SVM changes)
$wsThunderbird.Cells.Item(41, 6).Value = 10

# CV + tfidf + ngram(3) + POS
Set-ThunderbirdRow 45 72.36 69.27 64.48 83.26
Set-ThunderbirdRow 46 69.63 67.49 64.36 77.09
Set-ThunderbirdRow 47 5.47 48.63 3.76 10
Set-ThunderbirdRow 48 53.22 59.98 62.37 47.14
Set-ThunderbirdRow 49 51.53 62.47 70.4 41.55

# CV + tfidf + ngram(2) + POS
Set-ThunderbirdRow 51 71.56 69.54 65.77 79.27
Set-ThunderbirdRow 52 69.65 67.49 64.28 77.12
Set-ThunderbirdRow 53 5.47 48.63 3.76 10
Set-ThunderbirdRow 54 56.76 61.37 62.93 52.57
Set-ThunderbirdRow 55 55.06 63.67 70.19 46.02

# CV + tfidf + ngram(1) + POS
Set-ThunderbirdRow 57 69.22 69.02 67.44 71.99
Set-ThunderbirdRow 58 68.99 67.3 64.56 75.25
Set-ThunderbirdRow 59 5.47 48.63 3.76 10
Set-ThunderbirdRow 60 60.34 62.54 62.42 59.11
Set-ThunderbirdRow 61 60.55 65.87 69.84 54.18

# --- View / selection state ---
# Reviewer ends up back at the top of Thunderbird with F53 (the
# all-near-zero SVM row) selected, and that tab becomes the active one.
$wsLucene.Range("A1").Select()
$wsUbuntu.Range("F50").Select()

$wsThunderbird.Activate()
$wsThunderbird.Range("A1").Select()
$wsThunderbird.Range("F53").Select()
